$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = "42.167.38"
    $ws.Range("E2").Value = "  -1.09%  "
    $ws.Range("D3").Value = "2.241.83"
    $ws.Range("E3").Value = "  -1.10%  "
    $ws.Range("E4").Value = "  +0.21%  "
    $ws.Range("D5").Value = "246.36"
    $ws.Range("E5").Value = "  -1.70%  "
    $ws.Range("E6").Value = "  -0.38%  "
    $ws.Range("D7").Value = "74.16"
    $ws.Range("E7").Value = "  -2.74%  "
    $ws.Range("E8").Value = "  +0.13%  "
    $ws.Range("D9").Value = "0.615"
    $ws.Range("E9").Value = "  -4.37%  "
    $ws.Range("D10").Value = "41.79"
    $ws.Range("E10").Value = "  +4.29%  "
    $ws.Range("D11").Value = "0.0940"
    $ws.Range("E11").Value = "  -3.29%  "
    $ws.Range("D12").Value = "7.14"
    $ws.Range("E12").Value = "  -2.51%  "
    $ws.Range("E13").Value = "  -2.06%  "
    $ws.Range("D14").Value = "14.51"
    $ws.Range("E14").Value = "  -3.09%  "
    $ws.Range("D15").Value = "0.853"
    $ws.Range("E15").Value = "  -1.51%  "
    $ws.Range("D16").Value = "2.219.98"
    $ws.Range("E16").Value = "  -2.24%  "
    $ws.Range("D17").Value = "42.080.20"
    $ws.Range("E17").Value = "  -1.06%  "
    $ws.Range("D18").Value = "0.0₃0988"
    $ws.Range("E18").Value = "  -0.51%  "
    $ws.Range("E19").Value = "  -0.77%  "
    $ws.Range("D20").Value = "71.93"
    $ws.Range("E20").Value = "  -0.22%  "
    $ws.Range("D21").Value = "2.23"
    $ws.Range("E21").Value = "  +3.43%  "
    $ws.Range("D22").Value = "231.94"
    $ws.Range("E22").Value = "  -0.94%  "
    $ws.Range("D23").Value = "8.69"
    $ws.Range("E23").Value = "  +33.63%  "
    $ws.Range("E24").Value = "  +0.05%  "
    $ws.Range("D25").Value = "11.20"
    $ws.Range("E25").Value = "  -0.74%  "
    $ws.Range("D26").Value = "3.60"
    $ws.Range("E26").Value = "  -4.44%  "
    $ws.Range("E27").Value = "  -3.17%  "
    $ws.Range("B28").Value = "Toncoin"
    $ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    $ws.Range("D28").Value = "2.15"
    $ws.Range("E28").Value = "  +1.63%  "
    $ws.Range("B29").Value = "Monero"
    $ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    $ws.Range("D29").Value = "169.17"
    $ws.Range("E29").Value = "  +1.03%  "
    $ws.Range("D30").Value = "20.60"
    $ws.Range("E30").Value = "  -3.32%  "
    $ws.Range("D31").Value = "0.0821"
    $ws.Range("E31").Value = "  -4.23%  "
    $ws.Range("E32").Value = "  -3.90%  "
    $ws.Range("D33").Value = "30.39"
    $ws.Range("E33").Value = "  -3.53%  "
    $ws.Range("D34").Value = "0.125"
    $ws.Range("E34").Value = "  -1.61%  "
    $ws.Range("D35").Value = "5.15"
    $ws.Range("E35").Value = "  +8.80%  "
    $ws.Range("E36").Value = "  -0.58%  "
    $ws.Range("E37").Value = "  -0.73%  "
    $ws.Range("D38").Value = "13.78"
    $ws.Range("E38").Value = "  +0.69%  "
    $ws.Range("D39").Value = "2.18"
    $ws.Range("E39").Value = "  -3.53%  "
    $ws.Range("E40").Value = "  -1.38%  "
    $ws.Range("D41").Value = "62.05"
    $ws.Range("E41").Value = "  +0.35%  "
    $ws.Range("E42").Value = "  -3.03%  "
    $ws.Range("D43").Value = "107.35"
    $ws.Range("E43").Value = "  +1.02%  "
    $ws.Range("E44").Value = "  +1.68%  "
    $ws.Range("E45").Value = "  -2.59%  "
    $ws.Range("E46").Value = "  -0.14%  "
    $ws.Range("B47").Value = "ARBITRUM"
    $ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    $ws.Range("D47").Value = "1.12"
    $ws.Range("E47").Value = "  -2.36%  "
    $ws.Range("B48").Value = "FTXToken"
    $ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    $ws.Range("D48").Value = "4.34"
    $ws.Range("E48").Value = "  -8.27%  "
    $ws.Range("E49").Value = "  -0.63%  "
    $ws.Range("E50").Value = "  +1.21%  "
    $ws.Range("E51").Value = "  +0.14%  "
